$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("D1").Value = "State"
    $ws.Range("D2").Value = "STATE OF MndstcT TEXAS THL"
}
